$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (data rows start at row 2)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Column C ("Förändrad") holds a date serial number for every data row.
# The workbook update bumps this date from 2023-09-06 (45175) to 2023-09-08 (45177)
# for every data row, leaving everything else untouched.
$newDateSerial = 45177

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $newDateSerial
}
